# Fixed a bug in failure tracking

$p = $ppt.ActivePresentation

# --- Slide 12: title "Free Text Book (143 pages)" was split across three
# runs (leftover from an earlier edit); collapse it back into a single run
# while keeping the first run's formatting, without altering the text.
$s12 = $p.Slides.Item(12)
$title = $s12.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange

# First merge the 2nd+3rd runs' characters (a genuine structural change),
# then re-stamp the whole range so the engine collapses everything down to
# a single run using the first run's character formatting.
$titleRange.Characters(16, 11).Text = "(143 pages)"
$titleRange.Characters(1, $titleRange.Length).Text = "Free Text Book (143 pages)"

# --- Slide 7: "build your own interlock primitives" -> "... atomic ..."
$s7 = $p.Slides.Item(7)
$content = $s7.Shapes.Item(2)
$para = $content.TextFrame.TextRange.Paragraphs(4)
$para.Text = "build your own atomic primitives"
$para.Characters(1, $para.Length).Text = "build your own atomic primitives"
